$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.097.76'
$ws.Range("E2").Value = '  +4.02%  '

$ws.Range("D3").Value = '3.484.65'
$ws.Range("E3").Value = '  +3.66%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '410.12'
$ws.Range("E5").Value = '  -0.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.78'
$ws.Range("E6").Value = '  +18.79%  '

$ws.Range("D7").Value = '3.476.20'
$ws.Range("E7").Value = '  +3.63%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  +2.46%  '

$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("E10").Value = '  +9.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.128'
$ws.Range("E11").Value = '  +29.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.48'
$ws.Range("E12").Value = '  +8.47%  '

$ws.Range("E13").Value = '  -0.79%  '

$ws.Range("D14").Value = '4.027.53'
$ws.Range("E14").Value = '  +3.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.78'
$ws.Range("E15").Value = '  +2.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.33'
$ws.Range("E16").Value = '  +0.33%  '

$ws.Range("D17").Value = '3.523.32'
$ws.Range("E17").Value = '  +5.92%  '

$ws.Range("D18").Value = '62.926.96'
$ws.Range("E18").Value = '  +4.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.06'
$ws.Range("E19").Value = '  +0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.92'
$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000139'
$ws.Range("E21").Value = '  +25.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.38'
$ws.Range("E22").Value = '  -0.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '83.06'
$ws.Range("E23").Value = '  +9.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.28'
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '314.49'
$ws.Range("E25").Value = '  +3.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.20'
$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '30.59'
$ws.Range("E27").Value = '  +6.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.36'
$ws.Range("E28").Value = '  +4.80%  '

$ws.Range("E29").Value = '  +2.86%  '

$ws.Range("E30").Value = '  -0.61%  '

$ws.Range("E31").Value = '  -2.25%  '

$ws.Range("E32").Value = '  +4.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.88'
$ws.Range("E33").Value = '  +2.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.65'
$ws.Range("E34").Value = '  +0.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '43.20'
$ws.Range("E35").Value = '  +7.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0497'
$ws.Range("E37").Value = '  -2.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.58'
$ws.Range("E38").Value = '  +0.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.64'
$ws.Range("E39").Value = '  +6.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.03'
$ws.Range("E41").Value = '  -3.54%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.00'
$ws.Range("E42").Value = '  +3.62%  '

$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.126'
$ws.Range("E43").Value = '  +2.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '137.46'
$ws.Range("E44").Value = '  -0.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.67'
$ws.Range("E45").Value = '  +3.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.03'
$ws.Range("E46").Value = '  +0.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.289'
$ws.Range("E47").Value = '  -0.32%  '

$ws.Range("E48").Value = '  -1.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.48'
$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("D50").Value = '2.207.22'
$ws.Range("E50").Value = '  -0.11%  '

$ws.Range("D51").Value = '3.824.89'
$ws.Range("E51").Value = '  +3.70%  '
